# Auto-generated Excel COM-interop script
# Updates currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ
# figures across the per-job Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 280
$ws.Range("J12").Value = 270
$ws.Range("L12").Value = 270
$ws.Range("N12").Value = -610
$ws.Range("H62").Value = 5335.727
$ws.Range("I62").Value = 5742.2856
$ws.Range("J62").Value = 4624.25
$ws.Range("K62").Value = 5742.2856
$ws.Range("L62").Value = 4624.25
$ws.Range("M62").Value = -5118.2856
$ws.Range("N62").Value = -5872.25
$ws.Range("H65").Value = 5335.727
$ws.Range("I65").Value = 5742.2856
$ws.Range("J65").Value = 4624.25
$ws.Range("K65").Value = 28711.428
$ws.Range("L65").Value = 23121.25
$ws.Range("M65").Value = -25591.428
$ws.Range("N65").Value = -29361.25
$ws.Range("H132").Value = 7536.7754
$ws.Range("I132").Value = 7511.391
$ws.Range("J132").Value = 7926
$ws.Range("K132").Value = 22534.173
$ws.Range("L132").Value = 23778
$ws.Range("M132").Value = -20004.173
$ws.Range("N132").Value = -28838
$ws.Range("H137").Value = 36365520
$ws.Range("I137").Value = 20001908
$ws.Range("J137").Value = 200001630
$ws.Range("K137").Value = 60005724
$ws.Range("L137").Value = 600004890
$ws.Range("M137").Value = -60003174
$ws.Range("N137").Value = -600009990
$ws.Range("H138").Value = 5558369.5
$ws.Range("J138").Value = 7755251.5
$ws.Range("L138").Value = 23265754.5
$ws.Range("N138").Value = -23276034.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 100006280
$ws.Range("I61").Value = 200002560
$ws.Range("K61").Value = 200002560
$ws.Range("M61").Value = -200002348
$ws.Range("H122").Value = 2647.5386
$ws.Range("I122").Value = 1925.7858
$ws.Range("J122").Value = 3489.5833
$ws.Range("K122").Value = 5777.357400000001
$ws.Range("L122").Value = 10468.7499
$ws.Range("M122").Value = -3327.357400000001
$ws.Range("N122").Value = -15368.7499
$ws.Range("H132").Value = 23816068
$ws.Range("I132").Value = 7075.1577
$ws.Range("K132").Value = 21225.4731
$ws.Range("M132").Value = -18695.4731
$ws.Range("H136").Value = 100006280
$ws.Range("I136").Value = 200002560
$ws.Range("K136").Value = 600007680
$ws.Range("M136").Value = -600005130
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3604.0588
$ws.Range("I20").Value = 4178.091
$ws.Range("K20").Value = 4178.091
$ws.Range("M20").Value = -3931.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25645376
$ws.Range("I31").Value = 3337.8215
$ws.Range("J31").Value = 90916024
$ws.Range("K31").Value = 3337.8215
$ws.Range("L31").Value = 90916024
$ws.Range("M31").Value = -3042.8215
$ws.Range("N31").Value = -90916614
$ws.Range("H34").Value = 25645376
$ws.Range("I34").Value = 3337.8215
$ws.Range("J34").Value = 90916024
$ws.Range("K34").Value = 3337.8215
$ws.Range("L34").Value = 90916024
$ws.Range("M34").Value = -3135.8215
$ws.Range("N34").Value = -90916428
$ws.Range("H58").Value = 2045.862
$ws.Range("I58").Value = 1904.6428
$ws.Range("K58").Value = 1904.6428
$ws.Range("M58").Value = -1701.6428
$ws.Range("H105").Value = 11666.583
$ws.Range("J105").Value = 23392
$ws.Range("L105").Value = 23392
$ws.Range("N105").Value = -26886
$ws.Range("H107").Value = 1963.5264
$ws.Range("I107").Value = 1374.091
$ws.Range("K107").Value = 1374.091
$ws.Range("M107").Value = 545.9090000000001
$ws.Range("H136").Value = 2045.862
$ws.Range("I136").Value = 1904.6428
$ws.Range("K136").Value = 5713.928400000001
$ws.Range("M136").Value = -3163.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35863988
$ws.Range("I4").Value = 70928870
$ws.Range("K4").Value = 212786610
$ws.Range("M4").Value = -212786498
$ws.Range("H5").Value = 1648.84
$ws.Range("I5").Value = 639.875
$ws.Range("K5").Value = 1919.625
$ws.Range("M5").Value = -1807.625
$ws.Range("H26").Value = 642.2857
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 642.2857
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 1926.8571
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -2502.8571
$ws.Range("H37").Value = 199999
$ws.Range("J37").Value = 199999
$ws.Range("L37").Value = 599997
$ws.Range("N37").Value = -600221
$ws.Range("H113").Value = 4188.3335
$ws.Range("I113").Value = 3849
$ws.Range("K113").Value = 11547
$ws.Range("M113").Value = -9377
$ws.Range("H135").Value = 1648.84
$ws.Range("I135").Value = 639.875
$ws.Range("K135").Value = 5758.875
$ws.Range("M135").Value = -3223.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2365.55
$ws.Range("I102").Value = 2252.6428
$ws.Range("J102").Value = 2629
$ws.Range("K102").Value = 2252.6428
$ws.Range("L102").Value = 2629
$ws.Range("M102").Value = -630.6428000000001
$ws.Range("N102").Value = -5873
$ws.Range("H107").Value = 1188.8572
$ws.Range("J107").Value = 958
$ws.Range("L107").Value = 958
$ws.Range("N107").Value = -4798
$ws.Range("H122").Value = 2770.0667
$ws.Range("I122").Value = 2641
$ws.Range("K122").Value = 7923
$ws.Range("M122").Value = -5473
$ws.Range("H126").Value = 11116215
$ws.Range("I126").Value = 6672758
$ws.Range("K126").Value = 20018274
$ws.Range("M126").Value = -20015804
$ws.Range("H132").Value = 2764.2144
$ws.Range("I132").Value = 2484.9565
$ws.Range("K132").Value = 7454.869499999999
$ws.Range("M132").Value = -4924.869499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4578.7144
$ws.Range("I68").Value = 3875.5
$ws.Range("J68").Value = 5516.3335
$ws.Range("K68").Value = 3875.5
$ws.Range("L68").Value = 5516.3335
$ws.Range("M68").Value = -3126.5
$ws.Range("N68").Value = -7014.3335
$ws.Range("H71").Value = 4578.7144
$ws.Range("I71").Value = 3875.5
$ws.Range("J71").Value = 5516.3335
$ws.Range("K71").Value = 19377.5
$ws.Range("L71").Value = 27581.6675
$ws.Range("M71").Value = -15633.5
$ws.Range("N71").Value = -35069.6675
$ws.Range("H93").Value = 2466.9333
$ws.Range("I93").Value = 1267.5555
$ws.Range("J93").Value = 4266
$ws.Range("K93").Value = 1267.5555
$ws.Range("L93").Value = 4266
$ws.Range("M93").Value = -19.55549999999994
$ws.Range("N93").Value = -6762
$ws.Range("H100").Value = 2638.65
$ws.Range("I100").Value = 1720
$ws.Range("K100").Value = 1720
$ws.Range("M100").Value = -1179
$ws.Range("H132").Value = 74076330
$ws.Range("I132").Value = 2360.1667
$ws.Range("K132").Value = 7080.500100000001
$ws.Range("M132").Value = -4550.500100000001
$ws.Range("H136").Value = 3022.4055
$ws.Range("I136").Value = 3022.4055
$ws.Range("K136").Value = 9067.216499999999
$ws.Range("M136").Value = -6517.216499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3981.5151
$ws.Range("I132").Value = 4359.3584
$ws.Range("J132").Value = 2441.077
$ws.Range("K132").Value = 13078.0752
$ws.Range("L132").Value = 7323.231000000001
$ws.Range("M132").Value = -10548.0752
$ws.Range("N132").Value = -12383.231
$ws.Range("H136").Value = 2653.697
$ws.Range("I136").Value = 2487.9285
$ws.Range("J136").Value = 3582
$ws.Range("K136").Value = 7463.7855
$ws.Range("L136").Value = 10746
$ws.Range("M136").Value = -4913.7855
$ws.Range("N136").Value = -15846
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 85813.164
$ws.Range("J141").Value = 101484.5
$ws.Range("L141").Value = 101484.5
$ws.Range("N141").Value = -111844.5
